# Data update for past 3 days
# - Corrects E66 (tests_run) from 33304 -> 32304 and flags it in red (new cellXfs
#   style reusing the existing red font already present in the theme).
# - Appends three new daily rows (67, 68, 69) following the same formula
#   pattern as the existing rows, with F66:F69 sharing one formula group
#   (matches Excel's own behaviour when you fill a formula down a range).
# - Updates the sheet's selection to the new "last" data cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Correct the existing E66 value and mark it in red (Bad-style font).
# ---------------------------------------------------------------------------
$ws.Cells.Item(66, 5).Value = 32304
$ws.Cells.Item(66, 5).Font.Color = 255   # RGB(255,0,0) -> reuses the workbook's existing red font

# ---------------------------------------------------------------------------
# 2. Raw input values for the three new rows (67-69).
# ---------------------------------------------------------------------------
$ws.Cells.Item(67, 3).Value  = 52     # C67 new_cases
$ws.Cells.Item(67, 5).Value  = 32850  # E67 tests_run
$ws.Cells.Item(67, 11).Value = 168    # K67 hosp
$ws.Cells.Item(67, 12).Value = 191    # L67 death_count
$ws.Cells.Item(67, 15).Value = 14     # O67 ICU
$ws.Cells.Item(67, 17).Value = 1118   # Q67 active cases

$ws.Cells.Item(68, 3).Value  = 21     # C68
$ws.Cells.Item(68, 5).Value  = 33626  # E68
$ws.Cells.Item(68, 11).Value = 176    # K68
$ws.Cells.Item(68, 12).Value = 194    # L68
$ws.Cells.Item(68, 15).Value = 15     # O68
$ws.Cells.Item(68, 17).Value = 1094   # Q68

$ws.Cells.Item(69, 3).Value  = 28     # C69
$ws.Cells.Item(69, 5).Value  = 35059  # E69
$ws.Cells.Item(69, 11).Value = 181    # K69
$ws.Cells.Item(69, 12).Value = 194    # L69
$ws.Cells.Item(69, 15).Value = 15     # O69
$ws.Cells.Item(69, 17).Value = 1094   # Q69

# ---------------------------------------------------------------------------
# 3. Dates for column A (reuse the existing date-formatted style from A66
#    via copy/paste-formats, rather than minting a brand-new custom numFmt).
# ---------------------------------------------------------------------------
$ws.Cells.Item(67, 1).Value = 43974
$ws.Cells.Item(68, 1).Value = 43975
$ws.Cells.Item(69, 1).Value = 43976

$ws.Range("A66").Copy() | Out-Null
$ws.Range("A67:A69").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------------
# 4. Formulas that follow the existing per-row pattern (not shared - each
#    row's formula text has its own row-relative references, same as B65,
#    B66, etc. already stored in the sheet).
# ---------------------------------------------------------------------------
$ws.Cells.Item(67, 2).Formula  = "=B66+C67"
$ws.Cells.Item(67, 4).Formula  = "=((B67-B66)/B66)*100"
$ws.Cells.Item(67, 7).Formula  = "=((E67-E66)/E66)*100"
$ws.Cells.Item(67, 8).Formula  = "=B67"
$ws.Cells.Item(67, 9).Formula  = "=(H67/E67)*100"
$ws.Cells.Item(67, 10).Formula = "=(C67/F67)*100"
$ws.Cells.Item(67, 13).Formula = "=(L67/B67)*100"
$ws.Cells.Item(67, 14).Formula = "=((L67-L66)/L66)*100"
$ws.Cells.Item(67, 16).Formula = "=(O67/K67)*100"
$ws.Cells.Item(67, 18).Formula = "=Q67-K67"
$ws.Cells.Item(67, 19).Formula = "=(K67/Q67)*100"
$ws.Cells.Item(67, 20).Formula = "=(O67/Q67)*100"
$ws.Cells.Item(67, 21).Formula = "=K67-O67"

$ws.Cells.Item(68, 2).Formula  = "=B67+C68"
$ws.Cells.Item(68, 4).Formula  = "=((B68-B67)/B67)*100"
$ws.Cells.Item(68, 7).Formula  = "=((E68-E67)/E67)*100"
$ws.Cells.Item(68, 8).Formula  = "=B68"
$ws.Cells.Item(68, 9).Formula  = "=(H68/E68)*100"
$ws.Cells.Item(68, 10).Formula = "=(C68/F68)*100"
$ws.Cells.Item(68, 13).Formula = "=(L68/B68)*100"
$ws.Cells.Item(68, 14).Formula = "=((L68-L67)/L67)*100"
$ws.Cells.Item(68, 16).Formula = "=(O68/K68)*100"
$ws.Cells.Item(68, 18).Formula = "=Q68-K68"
$ws.Cells.Item(68, 19).Formula = "=(K68/Q68)*100"
$ws.Cells.Item(68, 20).Formula = "=(O68/Q68)*100"
$ws.Cells.Item(68, 21).Formula = "=K68-O68"

$ws.Cells.Item(69, 2).Formula  = "=B68+C69"
$ws.Cells.Item(69, 4).Formula  = "=((B69-B68)/B68)*100"
$ws.Cells.Item(69, 7).Formula  = "=((E69-E68)/E68)*100"
$ws.Cells.Item(69, 8).Formula  = "=B69"
$ws.Cells.Item(69, 9).Formula  = "=(H69/E69)*100"
$ws.Cells.Item(69, 10).Formula = "=(C69/F69)*100"
$ws.Cells.Item(69, 13).Formula = "=(L69/B69)*100"
$ws.Cells.Item(69, 14).Formula = "=((L69-L68)/L68)*100"
$ws.Cells.Item(69, 16).Formula = "=(O69/K69)*100"
$ws.Cells.Item(69, 18).Formula = "=Q69-K69"
$ws.Cells.Item(69, 19).Formula = "=(K69/Q69)*100"
$ws.Cells.Item(69, 20).Formula = "=(O69/Q69)*100"
$ws.Cells.Item(69, 21).Formula = "=K69-O69"

# ---------------------------------------------------------------------------
# 5. Column F (daily_tests) is filled as one contiguous range so it becomes
#    a shared formula group spanning F66:F69, exactly like Excel does when a
#    formula is dragged/filled down over previously-empty cells below it.
# ---------------------------------------------------------------------------
$ws.Range("F66:F69").Formula = "=E66-E65"

# ---------------------------------------------------------------------------
# 6. Update the active selection to reflect where the user ended up editing.
# ---------------------------------------------------------------------------
$ws.Range("F69").Select() | Out-Null
